$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels: kcal/mol -> kJ/mol
$ws.Range("G1").Value = "Relative Energy / kJ/mol"
$ws.Range("H1").Value = "Relative Energy + relative energy of [A]S_0 -> 2 [F]S_0 + O2 / kJ/mol"

# Update formulas to use the kJ/mol conversion factor (1 Hartree = 2625.5 kJ/mol)
$ws.Range("G2").Formula = "=((D2+D6)-(D4+D5))*2625.5"
$ws.Range("G3").Formula = "=((D3+2*D6)-(2*D4+2*D5))*2625.5"
$ws.Range("H3").Formula = "=G3+(82.728*4.184)"
$ws.Range("G7").Formula = "=((D7+2*D6)-(2*D5))*2625.5"

# Update selection to G1
$ws.Range("G1").Select()
